$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.711.55"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "2.260.09"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.99"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.26"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.49"
$ws.Range("E10").Value = "  +14.61%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.82"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "2.271.51"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "42.570.67"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("E18").Value = "  +4.34%  "
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.32"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.99"
$ws.Range("E21").Value = "  +52.01%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.01"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.93"
$ws.Range("E24").Value = "  +6.07%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.62"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +5.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.37"
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.75"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0831"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.15"
$ws.Range("E32").Value = "  -6.11%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.49"
$ws.Range("E33").Value = "  +12.43%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.120"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0318"
$ws.Range("E37").Value = "  +6.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.54"
$ws.Range("E38").Value = "  +7.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.21"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.85"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.56"
$ws.Range("E41").Value = "  +7.34%  "
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.11"
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.93"
$ws.Range("E44").Value = "  +1.85%  "
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.43"
$ws.Range("E47").Value = "  +8.45%  "
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.19"
$ws.Range("E51").Value = "  +0.81%  "
